$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.407.00"
$ws.Range('E2').Value = "'  -0.07%  "
$ws.Range('D3').Value = "'1.847.85"
$ws.Range('E3').Value = "'  -0.18%  "
$ws.Range('E4').Value = "'  +0.02%  "
$ws.Range('D5').Value = "'239.78"
$ws.Range('D6').Value = "'0.6318"
$ws.Range('E6').Value = "'  -0.21%  "
$ws.Range('D8').Value = "'0.07574"
$ws.Range('E8').Value = "'  -0.07%  "
$ws.Range('D9').Value = "'0.2935"
$ws.Range('E9').Value = "'  -0.90%  "
$ws.Range('D10').Value = "'24.59"
$ws.Range('E10').Value = "'  -0.10%  "
$ws.Range('E11').Value = "'  -0.15%  "
$ws.Range('D12').Value = "'1.885.98"
$ws.Range('E12').Value = "'  -4.99%  "
$ws.Range('D13').Value = "'5.005"
$ws.Range('E13').Value = "'  +0.10%  "
$ws.Range('D14').Value = "'0.6808"
$ws.Range('E14').Value = "'  -0.58%  "
$ws.Range('D15').Value = "'0.00001049"
$ws.Range('E15').Value = "'  +5.59%  "
$ws.Range('D16').Value = "'83.53"
$ws.Range('E16').Value = "'  +0.58%  "
$ws.Range('D17').Value = "'2.131.58"
$ws.Range('E17').Value = "'  -5.86%  "
$ws.Range('D18').Value = "'6.181"
$ws.Range('E18').Value = "'  -0.10%  "
$ws.Range('D19').Value = "'29.446.40"
$ws.Range('E19').Value = "'  -0.04%  "
$ws.Range('D20').Value = "'229.02"
$ws.Range('E20').Value = "'  -1.29%  "
$ws.Range('D21').Value = "'12.45"
$ws.Range('E21').Value = "'  -0.21%  "
$ws.Range('D22').Value = "'1.000"
$ws.Range('E22').Value = "'  +0.02%  "
$ws.Range('D23').Value = "'7.497"
$ws.Range('E23').Value = "'  -1.54%  "
$ws.Range('E24').Value = "'  +0.03%  "
$ws.Range('D25').Value = "'156.77"
$ws.Range('E25').Value = "'  +0.49%  "
$ws.Range('D26').Value = "'0.1396"
$ws.Range('E26').Value = "'  +0.51%  "
$ws.Range('D27').Value = "'8.342"
$ws.Range('E27').Value = "'  -0.77%  "
$ws.Range('E28').Value = "'  -0.47%  "
$ws.Range('D29').Value = "'1.461"
$ws.Range('E29').Value = "'  -0.65%  "
$ws.Range('D30').Value = "'1.301"
$ws.Range('E30').Value = "'  +3.40%  "
$ws.Range('D31').Value = "'0.05658"
$ws.Range('E31').Value = "'  -1.73%  "
$ws.Range('D32').Value = "'4.109"
$ws.Range('E32').Value = "'  -0.61%  "
$ws.Range('D33').Value = "'4.030"
$ws.Range('E33').Value = "'  +0.16%  "
$ws.Range('D34').Value = "'1.856"
$ws.Range('E34').Value = "'  -0.06%  "
$ws.Range('E35').Value = "'  -0.10%  "
$ws.Range('D36').Value = "'0.7105"
$ws.Range('E36').Value = "'  -0.87%  "
$ws.Range('E37').Value = "'  -0.16%  "
$ws.Range('D38').Value = "'1.249.69"
$ws.Range('E38').Value = "'  -0.47%  "
$ws.Range('D39').Value = "'0.01810"
$ws.Range('E39').Value = "'  +0.20%  "
$ws.Range('E40').Value = "'  -1.04%  "
$ws.Range('E41').Value = "'  +4.84%  "
$ws.Range('D42').Value = "'0.9028"
$ws.Range('E42').Value = "'  -0.26%  "
$ws.Range('E44').Value = "'  -0.03%  "
$ws.Range('E45').Value = "'  -1.59%  "
$ws.Range('B46').Value = "'Aptos"
$ws.Range('C46').Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range('D46').Value = "'7.095"
$ws.Range('E46').Value = "'  -0.76%  "
$ws.Range('B47').Value = "'TheSandbox"
$ws.Range('C47').Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range('D47').Value = "'0.4005"
$ws.Range('E47').Value = "'  -0.49%  "
$ws.Range('B48').Value = "'RenderToken"
$ws.Range('C48').Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range('D48').Value = "'1.678"
$ws.Range('E48').Value = "'  -0.48%  "
$ws.Range('B49').Value = "'EnergySwap"
$ws.Range('C49').Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range('D49').Value = "'8.952"
$ws.Range('E49').Value = "'  -2.60%  "
$ws.Range('B50').Value = "'Algorand"
$ws.Range('C50').Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range('D50').Value = "'0.1123"
$ws.Range('E50').Value = "'  -0.17%  "
$ws.Range('B51').Value = "'Cronos"
$ws.Range('C51').Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range('D51').Value = "'0.05704"
$ws.Range('E51').Value = "'  -0.71%  "
